$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

$ws.Cells.Item($row, 1).Value = "Partou"
$ws.Cells.Item($row, 2).Value = "KDV Partou De Groeskant 1"
$ws.Cells.Item($row, 3).Value = "KDV"

# Column D holds a plain "YYYY-MM-DD" text string in the source data (not a
# real date). Excel auto-parses date-looking input as a date serial when the
# cell is General-formatted, so format as Text first to force string storage,
# then clear the resulting format so the cell ends up unstyled (matching the
# rest of the sheet) while keeping the literal text value.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-10-14"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 1
